# fix: fix combination table
# Extends the "Combination" sheet with the missing rkus/rvip/rkir (1-4)
# Citta x Cetasika combination rows, mirroring the existing pattern used
# for every other Citta id already present in the table:
#   - one row with only the Citta id in column A
#   - one row with the Citta id in A and "karuna" in B
#   - one row with the Citta id in A and "mudita" in B

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combination")

$cittaIds = @("rkus1", "rkus2", "rkus3", "rkus4", "rvip1", "rvip2", "rvip3", "rvip4", "rkir1", "rkir2", "rkir3", "rkir4")

$row = 98
foreach ($cittaId in $cittaIds) {
    $ws.Cells.Item($row, 1).Value2 = $cittaId
    $row = $row + 1

    $ws.Cells.Item($row, 1).Value2 = $cittaId
    $ws.Cells.Item($row, 2).Value2 = "karuna"
    $row = $row + 1

    $ws.Cells.Item($row, 1).Value2 = $cittaId
    $ws.Cells.Item($row, 2).Value2 = "mudita"
    $row = $row + 1
}

$lastRow = $row - 1

# Extend the CetasikaID list validation to the newly added rows in column B.
$ws.Range("B98:B$lastRow").Validation.Add(3, 3, 1, "CetasikaID") | Out-Null

# Update the selection to match the new extent of the data, mirroring what
# Excel leaves behind after the rows were filled in.
$ws.Range("A1:D$lastRow").Select() | Out-Null
